$d = $word.ActiveDocument

# --- Professional Summary bullet: update the bold tech-stack sentence ---
# Original:  "...Java, Spring Boot, Spring Cloud, Spring MVC, ... Kubernetes, GCP cloud."
# Updated:   "...Java, Spring Boot, Spring MVC, Spring Cloud, Spring AI, ... Kubernetes, Mermaid & draw.io."
#
# Net effect: "Spring Cloud" moves after "Spring MVC" and gains a new
# "Spring AI" neighbour, and the trailing "GCP cloud" is swapped out for
# "Mermaid & draw.io" (with the closing period staying bold, same as the
# rest of the (now extended) bold run).

$wdReplaceAll = 2
$wdFindContinue = 1

$rng = $d.Content
$found = $rng.Find.Execute(
    "Spring Boot, Spring Cloud, Spring MVC, Micro-Service, Monolithic, hibernate, MySQL, MongoDB, Kafka, Redis, React JS, Redux, RTK, Micro-frontend & Docker, Kubernetes, GCP cloud.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Spring Boot, Spring MVC, Spring Cloud, Spring AI, Micro-Service, Monolithic, hibernate, MySQL, MongoDB, Kafka, Redis, React JS, Redux, RTK, Micro-frontend & Docker, Kubernetes, Mermaid & draw.io.",
    $wdReplaceAll)

if ($found) {
    # Keep the whole replaced span (stack text + trailing period) bold,
    # matching every run in the diff's target markup.
    $rng.Font.Bold = 1
}

Write-Host "Tech-stack sentence updated: $found"
